$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.135072333333333
$ws.Range("H2").Value = 21.405217
$ws.Range("I2").Value = 0.2435182897332695
$ws.Range("J2").Value = 0.2435182897332695
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 11.76385566666667
$ws.Range("N2").Value = 35.291567
$ws.Range("O2").Value = 0.3637526920523782
$ws.Range("P2").Value = 0.3637526920523783
$ws.Range("Q2").Value = 83.9359611005599
$ws.Range("R2").Value = 755.423649905039
$ws.Range("S2").Value = 0.08858043345446778
$ws.Range("T2").Value = 0.0885804334544678

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.135072333333333
$ws.Range("H3").Value = 21.405217
$ws.Range("I3").Value = 0.2435182897332695
$ws.Range("J3").Value = 0.2435182897332695
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.732885999999999
$ws.Range("N3").Value = 29.198658
$ws.Range("O3").Value = 0.3009526454809079
$ws.Range("P3").Value = 0.3009526454809079
$ws.Range("Q3").Value = 69.44484562208733
$ws.Range("R3").Value = 625.0036105987859
$ws.Range("S3").Value = 0.07328747351821366
$ws.Range("T3").Value = 0.07328747351821367

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.135072333333333
$ws.Range("H4").Value = 21.405217
$ws.Range("I4").Value = 0.2435182897332695
$ws.Range("J4").Value = 0.2435182897332695
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.575296333333333
$ws.Range("N4").Value = 16.725889
$ws.Range("O4").Value = 0.1723949279644982
$ws.Range("P4").Value = 0.1723949279644982
$ws.Range("Q4").Value = 39.78014261810144
$ws.Range("R4").Value = 358.021283562913
$ws.Range("S4").Value = 0.04198131801660478
$ws.Range("T4").Value = 0.04198131801660478

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 7.135072333333333
$ws.Range("H5").Value = 21.405217
$ws.Range("I5").Value = 0.2435182897332695
$ws.Range("J5").Value = 0.2435182897332695
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.268219333333334
$ws.Range("N5").Value = 15.804658
$ws.Range("O5").Value = 0.1628997345022157
$ws.Range("P5").Value = 0.1628997345022157
$ws.Range("Q5").Value = 37.58912601119845
$ws.Range("R5").Value = 338.302134100786
$ws.Range("S5").Value = 0.03966906474398323
$ws.Range("T5").Value = 0.03966906474398324

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.272738
$ws.Range("H6").Value = 33.818214
$ws.Range("I6").Value = 0.3847358162785133
$ws.Range("J6").Value = 0.3847358162785133
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 11.76385566666667
$ws.Range("N6").Value = 35.291567
$ws.Range("O6").Value = 0.3637526920523782
$ws.Range("P6").Value = 0.3637526920523783
$ws.Range("Q6").Value = 132.6108628001487
$ws.Range("R6").Value = 1193.497765201338
$ws.Range("S6").Value = 0.1399486889002784
$ws.Range("T6").Value = 0.1399486889002784

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.272738
$ws.Range("H7").Value = 33.818214
$ws.Range("I7").Value = 0.3847358162785133
$ws.Range("J7").Value = 0.3847358162785133
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.732885999999999
$ws.Range("N7").Value = 29.198658
$ws.Range("O7").Value = 0.3009526454809079
$ws.Range("P7").Value = 0.3009526454809079
$ws.Range("Q7").Value = 109.716273861868
$ws.Range("R7").Value = 987.4464647568119
$ws.Range("S7").Value = 0.1157872617202751
$ws.Range("T7").Value = 0.1157872617202751

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.272738
$ws.Range("H8").Value = 33.818214
$ws.Range("I8").Value = 0.3847358162785133
$ws.Range("J8").Value = 0.3847358162785133
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.575296333333333
$ws.Range("N8").Value = 16.725889
$ws.Range("O8").Value = 0.1723949279644982
$ws.Range("P8").Value = 0.1723949279644982
$ws.Range("Q8").Value = 62.84885483802732
$ws.Range("R8").Value = 565.6396935422459
$ws.Range("S8").Value = 0.06632650333269668
$ws.Range("T8").Value = 0.06632650333269668

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.272738
$ws.Range("H9").Value = 33.818214
$ws.Range("I9").Value = 0.3847358162785133
$ws.Range("J9").Value = 0.3847358162785133
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.268219333333334
$ws.Range("N9").Value = 15.804658
$ws.Range("O9").Value = 0.1628997345022157
$ws.Range("P9").Value = 0.1628997345022157
$ws.Range("Q9").Value = 59.38725627120133
$ws.Range("R9").Value = 534.485306440812
$ws.Range("S9").Value = 0.06267336232526304
$ws.Range("T9").Value = 0.06267336232526305

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.311962000000001
$ws.Range("H10").Value = 18.935886
$ws.Range("I10").Value = 0.2154257335164676
$ws.Range("J10").Value = 0.2154257335164676
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 11.76385566666667
$ws.Range("N10").Value = 35.291567
$ws.Range("O10").Value = 0.3637526920523782
$ws.Range("P10").Value = 0.3637526920523783
$ws.Range("Q10").Value = 74.25300994148468
$ws.Range("R10").Value = 668.2770894733621
$ws.Range("S10").Value = 0.07836169050397333
$ws.Range("T10").Value = 0.07836169050397333

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 6.311962000000001
$ws.Range("H11").Value = 18.935886
$ws.Range("I11").Value = 0.2154257335164676
$ws.Range("J11").Value = 0.2154257335164676
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 9.732885999999999
$ws.Range("N11").Value = 29.198658
$ws.Range("O11").Value = 0.3009526454809079
$ws.Range("P11").Value = 0.3009526454809079
$ws.Range("Q11").Value = 61.433606582332
$ws.Range("R11").Value = 552.902459240988
$ws.Range("S11").Value = 0.064832944406446
$ws.Range("T11").Value = 0.06483294440644602

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 6.311962000000001
$ws.Range("H12").Value = 18.935886
$ws.Range("I12").Value = 0.2154257335164676
$ws.Range("J12").Value = 0.2154257335164676
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 5.575296333333333
$ws.Range("N12").Value = 16.725889
$ws.Range("O12").Value = 0.1723949279644982
$ws.Range("P12").Value = 0.1723949279644982
$ws.Range("Q12").Value = 35.19105859473934
$ws.Range("R12").Value = 316.719527352654
$ws.Range("S12").Value = 0.0371383038112706
$ws.Range("T12").Value = 0.0371383038112706

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 6.311962000000001
$ws.Range("H13").Value = 18.935886
$ws.Range("I13").Value = 0.2154257335164676
$ws.Range("J13").Value = 0.2154257335164676
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 5.268219333333334
$ws.Range("N13").Value = 15.804658
$ws.Range("O13").Value = 0.1628997345022157
$ws.Range("P13").Value = 0.1628997345022157
$ws.Range("Q13").Value = 33.25280023966534
$ws.Range("R13").Value = 299.2752021569881
$ws.Range("S13").Value = 0.03509279479477764
$ws.Range("T13").Value = 0.03509279479477764

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.580172
$ws.Range("H14").Value = 13.740516
$ws.Range("I14").Value = 0.1563201604717497
$ws.Range("J14").Value = 0.1563201604717497
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 11.76385566666667
$ws.Range("N14").Value = 35.291567
$ws.Range("O14").Value = 0.3637526920523782
$ws.Range("P14").Value = 0.3637526920523783
$ws.Range("Q14").Value = 53.880482336508
$ws.Range("R14").Value = 484.924341028572
$ws.Range("S14").Value = 0.05686187919365872
$ws.Range("T14").Value = 0.05686187919365871

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.580172
$ws.Range("H15").Value = 13.740516
$ws.Range("I15").Value = 0.1563201604717497
$ws.Range("J15").Value = 0.1563201604717497
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 9.732885999999999
$ws.Range("N15").Value = 29.198658
$ws.Range("O15").Value = 0.3009526454809079
$ws.Range("P15").Value = 0.3009526454809079
$ws.Range("Q15").Value = 44.578291936392
$ws.Range("R15").Value = 401.204627427528
$ws.Range("S15").Value = 0.04704496583597312
$ws.Range("T15").Value = 0.04704496583597312

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.580172
$ws.Range("H16").Value = 13.740516
$ws.Range("I16").Value = 0.1563201604717497
$ws.Range("J16").Value = 0.1563201604717497
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 5.575296333333333
$ws.Range("N16").Value = 16.725889
$ws.Range("O16").Value = 0.1723949279644982
$ws.Range("P16").Value = 0.1723949279644982
$ws.Range("Q16").Value = 25.535816157636
$ws.Range("R16").Value = 229.822345418724
$ws.Range("S16").Value = 0.02694880280392608
$ws.Range("T16").Value = 0.02694880280392608

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.580172
$ws.Range("H17").Value = 13.740516
$ws.Range("I17").Value = 0.1563201604717497
$ws.Range("J17").Value = 0.1563201604717497
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 5.268219333333334
$ws.Range("N17").Value = 15.804658
$ws.Range("O17").Value = 0.1628997345022157
$ws.Range("P17").Value = 0.1628997345022157
$ws.Range("Q17").Value = 24.129350680392
$ws.Range("R17").Value = 217.164156123528
$ws.Range("S17").Value = 0.02546451263819178
$ws.Range("T17").Value = 0.02546451263819178
